# Atualizando o arquivo XLSX
# - S2 (Odd_Over05_HT for row 2) changes from 1.47 to 1.5
# - A new match row (row 4) is appended for URUGUAY - PRIMERA DIVISION:
#   Cerro Largo vs Liverpool M.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing odd value in S2
$ws.Range("S2").Value = 1.5

# Columns A..BD (56 columns) of the new row 4, in header order:
# Id, Date, Time, League, Home, Away,
# Odd_H_FT, Odd_D_FT, Odd_A_FT, Odd_H_HT, Odd_D_HT, Odd_A_HT,
# Odd_Over05_FT, Odd_Under05_FT, Odd_Over15_FT, Odd_Under15_FT,
# Odd_Over25_FT, Odd_Under25_FT, Odd_Over05_HT, Odd_Under05_HT,
# Odd_BTTS_Yes, Odd_BTTS_No,
# Odd_CS_1-0, Odd_CS_2-0, Odd_CS_2-1, Odd_CS_3-0, Odd_CS_3-1, Odd_CS_3-2,
# Odd_CS_0-0, Odd_CS_1-1, Odd_CS_2-2, Odd_CS_3-3,
# Odd_CS_0-1, Odd_CS_0-2, Odd_CS_1-2, Odd_CS_0-3, Odd_CS_1-3, Odd_CS_2-3,
# Odd_CS_1-0_HT, Odd_CS_2-0_HT, Odd_CS_2-1_HT, Odd_CS_3-0_HT, Odd_CS_3-1_HT,
# Odd_CS_0-0_HT, Odd_CS_1-1_HT, Odd_CS_2-2_HT,
# Odd_CS_0-1_HT, Odd_CS_0-2_HT, Odd_CS_1-2_HT, Odd_CS_0-3_HT, Odd_CS_1-3_HT,
# Odd_CS_2-3_HT, Odd_CS_4-4, Odd_CS_3-2_HT, Odd_CS_3-3_HT, Odd_CS_4-4_HT
$row4Values = @(
    'QXb6IStp','13/11/2024','19:00','URUGUAY - PRIMERA DIVISION','Cerro Largo','Liverpool M.',
    2.7,3,2.88,3.5,1.95,3.6,
    1.1,7,1.44,2.63,
    2.4,1.53,1.53,2.38,
    2,1.73,
    8,13,11,29,26,34,
    8,6,15,51,
    7.5,12,10,26,23,34,
    4.75,17,29,51,81,
    1.83,8.5,51,
    4.5,15,26,51,81,
    201,301,201,51,51
)

for ($i = 0; $i -lt $row4Values.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4Values[$i]
}
